$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 142
$ws1.Range("F4").Value = 1768
$ws1.Range("F6").Value = 1038
$ws1.Range("F7").Value = 2193
$ws1.Range("F8").Value = 2108
$ws1.Range("F9").Value = 1104
$ws1.Range("F10").Value = 604
$ws1.Range("F11").Value = 20
$ws1.Range("F12").Value = 1669
$ws1.Range("F18").Value = 1579
$ws1.Range("F19").Value = 629
$ws1.Range("F22").Value = 12221
$ws1.Range("F23").Value = 12262
$ws1.Range("F25").Value = 700
$ws1.Range("F27").Value = 30
$ws1.Range("F29").Value = 362
$ws1.Range("F32").Value = 574

# --- Sheet "演出" (sheet2) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 40
$ws2.Range("G5").Value = 180
$ws2.Range("F6").Value = 40

# --- Sheet "全部类型" (sheet4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 142
$ws4.Range("F5").Value = 1768
$ws4.Range("F7").Value = 1038
$ws4.Range("F8").Value = 2193
$ws4.Range("F9").Value = 2108
$ws4.Range("F10").Value = 1104
$ws4.Range("F11").Value = 604
$ws4.Range("F12").Value = 20
$ws4.Range("F13").Value = 1669
$ws4.Range("F20").Value = 40
$ws4.Range("F22").Value = 1579
$ws4.Range("F23").Value = 629
$ws4.Range("F26").Value = 12221
$ws4.Range("F27").Value = 12262
$ws4.Range("F29").Value = 700
$ws4.Range("F31").Value = 30
$ws4.Range("F33").Value = 362
$ws4.Range("G35").Value = 180
$ws4.Range("F36").Value = 40
$ws4.Range("F38").Value = 574
